$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("threshold_b")
$ws.Activate()
Write-Host "Active sheet now:" $wb.ActiveSheet.Name
